$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The workbook has a "totals" row as the last row (row 30), summing the
# hours of all the logged days above it. A new logged day is being added,
# so that totals row needs to move down to row 31 (with an updated SUM
# range), and the new entry takes its old place at row 30.

# Step 1: create the new totals row at 31, using the old totals row (30) as
# a formatting template (so style/border/fill match the original totals row).
$ws.Range("B30:D30").Copy()
$ws.Range("B31:D31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(31, 2).Value = "Yht"
$ws.Cells.Item(31, 4).Value = $null
$ws.Rows.Item(31).RowHeight = 18.75

# Step 2: turn the old totals row (30) into a normal data entry row, using
# the row above it (29, a typical data row) as the formatting template.
$ws.Range("B29:D29").Copy()
$ws.Range("B30:D30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$newDate = Get-Date -Year 2024 -Month 2 -Day 29 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(30, 2).Value = $newDate
$ws.Cells.Item(30, 3).Value = 4
$ws.Cells.Item(30, 4).Value = "Paransin systeminformation tyylitelyä ja lataamista saadakseen käytäjä kokemuksen paremaksi."
$ws.Rows.Item(30).RowHeight = 37.5

# Step 3: now that the new entry is in place, (re)point the totals formula
# at the expanded range.
$ws.Cells.Item(31, 3).Formula = "=SUM(C6:C30)"

# Update the sheet view to match where the user ended up scrolled/selected.
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("F29").Select()

$wb.Save()
